$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '24.254.03'
Set-TextCell 'E2' '  +11.39%  '
Set-TextCell 'D3' '1.683.77'
Set-TextCell 'E3' '  +7.20%  '
Set-TextCell 'E4' '  +0.42%  '
Set-TextCell 'D5' '309.84'
Set-TextCell 'E5' '  +8.65%  '
Set-TextCell 'D6' '0.9983'
Set-TextCell 'E6' '  +2.00%  '
Set-TextCell 'D7' '0.3747'
Set-TextCell 'E7' '  +1.80%  '
Set-TextCell 'E8' '  +5.88%  '
Set-TextCell 'E9' '  +15.90%  '
Set-TextCell 'E10' '  +6.17%  '
Set-TextCell 'D11' '0.07340'
Set-TextCell 'E11' '  +4.27%  '
Set-TextCell 'D12' '0.9991'
Set-TextCell 'E12' '  +0.46%  '
Set-TextCell 'D13' '20.60'
Set-TextCell 'E13' '  +2.80%  '
Set-TextCell 'D14' '6.156'
Set-TextCell 'E14' '  +6.04%  '
Set-TextCell 'D15' '6.807'
Set-TextCell 'E15' '  +5.11%  '
Set-TextCell 'D16' '1.684.88'
Set-TextCell 'E16' '  +7.72%  '
Set-TextCell 'D17' '0.00001114'
Set-TextCell 'E17' '  +4.51%  '
Set-TextCell 'E18' '  +2.23%  '
Set-TextCell 'D19' '0.06744'
Set-TextCell 'E19' '  +9.25%  '
Set-TextCell 'D20' '82.61'
Set-TextCell 'E20' '  +11.68%  '
Set-TextCell 'E21' '  +3.25%  '
Set-TextCell 'D22' '6.141'
Set-TextCell 'E22' '  +5.40%  '
Set-TextCell 'D23' '12.10'
Set-TextCell 'E23' '  +4.52%  '
Set-TextCell 'D24' '24.217.11'
Set-TextCell 'E24' '  +11.24%  '
Set-TextCell 'D25' '2.416'
Set-TextCell 'E25' '  +3.27%  '
Set-TextCell 'D26' '2.701'
Set-TextCell 'E26' '  +12.46%  '
Set-TextCell 'D27' '3.366'
Set-TextCell 'E27' '  -8.75%  '
Set-TextCell 'D28' '153.08'
Set-TextCell 'E28' '  +3.19%  '
Set-TextCell 'D29' '19.72'
Set-TextCell 'E29' '  +8.18%  '
Set-TextCell 'D30' '1.867.78'
Set-TextCell 'E30' '  +7.45%  '
Set-TextCell 'D31' '127.28'
Set-TextCell 'E31' '  +6.20%  '
Set-TextCell 'D32' '6.500'
Set-TextCell 'E32' '  +20.74%  '
Set-TextCell 'D33' '4.079'
Set-TextCell 'E33' '  +0.16%  '
Set-TextCell 'D34' '0.9991'
Set-TextCell 'E34' '  +9.99%  '
Set-TextCell 'D35' '1.792'
Set-TextCell 'E35' '  +12.91%  '
Set-TextCell 'D36' '0.08531'
Set-TextCell 'E36' '  +4.35%  '
Set-TextCell 'E37' '  +9.48%  '
Set-TextCell 'D38' '0.06521'
Set-TextCell 'D39' '5.418'
Set-TextCell 'E39' '  +5.88%  '
Set-TextCell 'D40' '9.005'
Set-TextCell 'E40' '  +11.09%  '
Set-TextCell 'D41' '0.02366'
Set-TextCell 'E41' '  +9.58%  '
Set-TextCell 'D42' '1.283'
Set-TextCell 'E42' '  +4.31%  '
Set-TextCell 'D43' '0.2154'
Set-TextCell 'E43' '  +7.71%  '
Set-TextCell 'D44' '0.6260'
Set-TextCell 'E44' '  +9.46%  '
Set-TextCell 'D45' '0.9976'
Set-TextCell 'E45' '  +2.35%  '
Set-TextCell 'B46' 'PancakeSwap'
Set-TextCell 'C46' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D46' '3.816'
Set-TextCell 'E46' '  +5.78%  '
Set-TextCell 'B47' 'EnergySwap'
Set-TextCell 'C47' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D47' '13.26'
Set-TextCell 'E47' '  +3.97%  '
Set-TextCell 'D48' '0.5998'
Set-TextCell 'E48' '  +6.90%  '
Set-TextCell 'D49' '128.36'
Set-TextCell 'E49' '  +3.22%  '
Set-TextCell 'D50' '2.048'
Set-TextCell 'E50' '  +7.00%  '
Set-TextCell 'D51' '0.07196'
Set-TextCell 'E51' '  +6.85%  '
